$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comments")
$ws.Range("G2:G72").Validation.Add(3, 1, 1, "='(list)'!`$A`$1:`$A`$22")
Write-Host "added"
